$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 4 data rows (rows 2-5, the oldest observations).
# This shifts all subsequent rows up by 4, which also drops the last
# 4 rows off the bottom of the data range (old rows 40-43).
$ws.Range("A2:B5").EntireRow.Delete()
